$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Ludmila
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Ludmila"
$ws.Range("C4").Value = "<Cachorro.Cachorro object at 0x0000019680FF9DC0>"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "31995389240"
$ws.Range("E4").Value = 0

# Row 5 - Laura
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Laura"
$ws.Range("C5").Value = 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "31295859949"
$ws.Range("E5").Value = 0

# Row 6 - Paulo Nogueira
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "Paulo Nogueira"
$ws.Range("C6").Value = 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "31907838920"
$ws.Range("E6").Value = 0

# Row 7 - Camila Figueredo
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Camila Figueredo"
$ws.Range("C7").Value = 7
$ws.Range("D7").Value = "Rua Lameda dos Perdizes 23"
$ws.Range("E7").Value = 0

# Row 8 - Lilian Campos
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Lilian Campos"
$ws.Range("C8").Value = 8
$ws.Range("D8").Value = "Rua Geraldo Luiz 90"
$ws.Range("E8").Value = 0

# Row 9 - Lucas Gontijo
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Lucas Gontijo"
$ws.Range("C9").Value = 9
$ws.Range("D9").Value = "Rua Claudio Bandeira 45"
$ws.Range("E9").Value = 0

# Row 10 - Saulo Nunes
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Saulo Nunes"
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = "Rua Tajubá 789"
$ws.Range("E10").Value = 0

# Row 11 - Fausto Silva
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Fausto Silva"
$ws.Range("C11").Value = 11
$ws.Range("D11").Value = "Rua Bernardo Monteiro 638"
$ws.Range("E11").Value = 0

# Row 12 - Otavio Marquez
$ws.Range("A12").Value = 12
$ws.Range("B12").Value = "Otavio Marquez"
$ws.Range("C12").Value = 12
$ws.Range("D12").Value = "Rua Carlos Luz 34"
$ws.Range("E12").Value = 0
